$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix product description typo: spf50+ -> spf30+
$ws.Range("E28").Value = "Nivea Sun® Kids Spray TRIGGER spf30+ 300ml"

# Update "ΝΕΑ ΤΙΜΗ" (New Price) column G values
$ws.Range("G2").Value = 16.2
$ws.Range("G3").Value = 14.9
$ws.Range("G4").Value = 13.9
$ws.Range("G5").Value = 13.9
$ws.Range("G6").Value = 13.9
$ws.Range("G7").Value = 14.8
$ws.Range("G8").Value = 14.2
$ws.Range("G9").Value = 15.2
$ws.Range("G10").Value = 15.45
$ws.Range("G11").Value = 9.9
$ws.Range("G12").Value = 9.9
$ws.Range("G17").Value = 15.98
$ws.Range("G20").Value = 14.95
$ws.Range("G21").Value = 16.4
$ws.Range("G26").Value = 13.55
$ws.Range("G28").Value = 13.4
$ws.Range("G30").Value = 15.2
$ws.Range("G32").Value = 9.9
$ws.Range("G36").Value = 14.95
$ws.Range("G37").Value = 10.9
$ws.Range("G46").Value = 17.9
